$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Status (B2, C2) and Latest Handoff Date (D2)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-29-21 02:29:24"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (E2)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-21 02:29:20"

# de-de sheet: Status (C2) and Latest Handoff Datetime (E2)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-21 02:29:24"
